$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-record data between row 3 and row 5
# (columns D, M, N, O, P, Q, R, S, T) while leaving the shared
# identification columns (A, B, C, E-L) untouched.

$row3_D = $ws.Range("D3").Value2
$row3_M = $ws.Range("M3").Value2
$row3_N = $ws.Range("N3").Value2
$row3_O = $ws.Range("O3").Value2
$row3_P = $ws.Range("P3").Value2
$row3_Q = $ws.Range("Q3").Value2
$row3_R = $ws.Range("R3").Value2
$row3_S = $ws.Range("S3").Value2
$row3_T = $ws.Range("T3").Value2

$row5_D = $ws.Range("D5").Value2
$row5_M = $ws.Range("M5").Value2
$row5_N = $ws.Range("N5").Value2
$row5_O = $ws.Range("O5").Value2
$row5_P = $ws.Range("P5").Value2
$row5_Q = $ws.Range("Q5").Value2
$row5_R = $ws.Range("R5").Value2
$row5_S = $ws.Range("S5").Value2
$row5_T = $ws.Range("T5").Value2

$ws.Range("D3").Value = $row5_D
$ws.Range("M3").Value = $row5_M
$ws.Range("N3").Value = $row5_N
$ws.Range("O3").Value = $row5_O
$ws.Range("P3").Value = $row5_P
$ws.Range("Q3").Value = $row5_Q
$ws.Range("R3").Value = $row5_R
$ws.Range("S3").Value = $row5_S
$ws.Range("T3").Value = $row5_T

$ws.Range("D5").Value = $row3_D
$ws.Range("M5").Value = $row3_M
$ws.Range("N5").Value = $row3_N
$ws.Range("O5").Value = $row3_O
$ws.Range("P5").Value = $row3_P
$ws.Range("Q5").Value = $row3_Q
$ws.Range("R5").Value = $row3_R
$ws.Range("S5").Value = $row3_S
$ws.Range("T5").Value = $row3_T
